$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 2.2
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 2.88
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("S2").Value = 1.5
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 1.83
$ws.Range("AA2").Value = 19
$ws.Range("AC2").Value = 8
$ws.Range("AG2").Value = 301
$ws.Range("AH2").Value = 9.5
$ws.Range("AK2").Value = 41
$ws.Range("AN2").Value = 4
$ws.Range("AP2").Value = 23
$ws.Range("AV2").Value = 51

# Row 3
$ws.Range("G3").Value = 5.7
$ws.Range("H3").Value = 4.4
$ws.Range("I3").Value = 1.42
$ws.Range("J3").Value = 5.4
$ws.Range("K3").Value = 2.45
$ws.Range("L3").Value = 1.87
$ws.Range("P3").Value = 4.35
$ws.Range("Q3").Value = 1.52
$ws.Range("R3").Value = 2.22
$ws.Range("U3").Value = 1.73
$ws.Range("V3").Value = 2.06
$ws.Range("W3").Value = 15.5
$ws.Range("X3").Value = 30
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 80
$ws.Range("AA3").Value = 40
$ws.Range("AB3").Value = 37
$ws.Range("AC3").Value = 15
$ws.Range("AD3").Value = 7.8
$ws.Range("AH3").Value = 7.3
$ws.Range("AI3").Value = 6.6
$ws.Range("AK3").Value = 8.5
$ws.Range("AL3").Value = 9
$ws.Range("AN3").Value = 7.5
$ws.Range("AO3").Value = 30
$ws.Range("AP3").Value = 32
$ws.Range("AT3").Value = 3.35
$ws.Range("AU3").Value = 7.6
$ws.Range("AW3").Value = 3.4
$ws.Range("AX3").Value = 6.4
$ws.Range("AZ3").Value = 17
$ws.Range("BB3").Value = 150

# Row 5
$ws.Range("G5").Value = 2.35
$ws.Range("I5").Value = 2.8
$ws.Range("J5").Value = 2.92
$ws.Range("L5").Value = 3.25
$ws.Range("W5").Value = 7.9
$ws.Range("X5").Value = 10.75
$ws.Range("Y5").Value = 7.7
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 15
$ws.Range("AB5").Value = 18.5
$ws.Range("AD5").Value = 5.5
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 9.75
$ws.Range("AI5").Value = 14.5
$ws.Range("AJ5").Value = 8.5
$ws.Range("AK5").Value = 29
$ws.Range("AL5").Value = 17
$ws.Range("AN5").Value = 4.5
$ws.Range("AO5").Value = 12.5
$ws.Range("AP5").Value = 18
$ws.Range("AQ5").Value = 50
$ws.Range("AR5").Value = 75
$ws.Range("AS5").Value = 175
$ws.Range("AU5").Value = 6.1
$ws.Range("AW5").Value = 5
$ws.Range("AX5").Value = 14.5
$ws.Range("AY5").Value = 18
$ws.Range("AZ5").Value = 60
$ws.Range("BA5").Value = 75

# Row 6
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 8

# Row 7
$ws.Range("G7").Value = 2.9
$ws.Range("I7").Value = 2.35
$ws.Range("J7").Value = 3.4
$ws.Range("M7").Value = 1.04
$ws.Range("O7").Value = 1.22
$ws.Range("Y7").Value = 11
$ws.Range("AB7").Value = 29
$ws.Range("AI7").Value = 12
$ws.Range("AL7").Value = 17
$ws.Range("AM7").Value = 23
$ws.Range("AP7").Value = 23

# Row 8
$ws.Range("M8").Value = 1.04
$ws.Range("O8").Value = 1.22

# Row 9
$ws.Range("J9").Value = 2.4
$ws.Range("M9").Value = 1.04
$ws.Range("O9").Value = 1.25
$ws.Range("AA9").Value = 15
$ws.Range("AC9").Value = 12
$ws.Range("AY9").Value = 29

# Row 10
$ws.Range("M10").Value = 1.03
$ws.Range("O10").Value = 1.2
$ws.Range("P10").Value = 4.33
$ws.Range("Q10").Value = 1.65
$ws.Range("R10").Value = 2.2
$ws.Range("AH10").Value = 9.5
$ws.Range("AL10").Value = 15
$ws.Range("AU10").Value = 7.5

# Row 11
$ws.Range("O11").Value = 1.18
$ws.Range("P11").Value = 4.5
$ws.Range("Q11").Value = 1.65
$ws.Range("R11").Value = 2.2

# Row 14
$ws.Range("Q14").Value = 2.4
$ws.Range("R14").Value = 1.53

# Row 15
$ws.Range("G15").Value = 1.91
$ws.Range("H15").Value = 2.9
$ws.Range("I15").Value = 4.1
$ws.Range("S15").Value = 1.62
$ws.Range("T15").Value = 2.2
$ws.Range("X15").Value = 8
$ws.Range("Y15").Value = 10
$ws.Range("Z15").Value = 17
$ws.Range("AH15").Value = 8.5
$ws.Range("AI15").Value = 19
$ws.Range("AO15").Value = 12
$ws.Range("AT15").Value = 2.2
$ws.Range("AX15").Value = 26

# Row 16
$ws.Range("M16").Value = 1.07
$ws.Range("N16").Value = 9

# Row 17
$ws.Range("G17").Value = 1.85
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 2.5
$ws.Range("K17").Value = 2.2
$ws.Range("Z17").Value = 15
$ws.Range("AC17").Value = 10
$ws.Range("AH17").Value = 11
$ws.Range("AI17").Value = 21
$ws.Range("AN17").Value = 3.75
$ws.Range("AX17").Value = 23

# Row 18
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 3.25
$ws.Range("I18").Value = 3.6
$ws.Range("J18").Value = 2.75
$ws.Range("Q18").Value = 2.1
$ws.Range("R18").Value = 1.7
$ws.Range("AC18").Value = 8.5
$ws.Range("AE18").Value = 17
$ws.Range("AH18").Value = 10
$ws.Range("AI18").Value = 17
$ws.Range("AP18").Value = 23
$ws.Range("AU18").Value = 8.5
$ws.Range("BB18").Value = 251

# Row 19
$ws.Range("G19").Value = 2.38
$ws.Range("I19").Value = 3.3
$ws.Range("M19").Value = 1.1
$ws.Range("N19").Value = 7
$ws.Range("U19").Value = 2
$ws.Range("V19").Value = 1.73
$ws.Range("Z19").Value = 23
$ws.Range("AD19").Value = 5.5
$ws.Range("AQ19").Value = 51
$ws.Range("AR19").Value = 81
$ws.Range("AS19").Value = 251
